$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A7")
$rng.Interior.Pattern = -4105
